$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header summary fields -------------------------------------------------
# Valor Mora total
$ws.Range("E11").Value = 216030

# Cant. Trabajadores / Cant. Periodos
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 6

# --- Remove the now-obsolete detail rows ------------------------------------
# Old table had 17 data rows (16-32); the new statement only needs 6 data
# rows, so drop the 11 rows in the middle (21-31). Row 32 (with its closing
# bottom-border styling) slides up to become the new row 21.
$ws.Rows("21:31").Delete()

# --- Rewrite the 6 remaining detail rows with the new worker / periods -----
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "45761661"
$ws.Range("D16").Value = "GLORIA ENIT DIAZ AVILA"
$ws.Range("E16").Value = "2107"
$ws.Range("F16").Value = 36341
$ws.Range("G16").Value = 781242

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "45761661"
$ws.Range("D17").Value = "GLORIA ENIT DIAZ AVILA"
$ws.Range("E17").Value = "2110"
$ws.Range("F17").Value = 36341
$ws.Range("G17").Value = 781242

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "45761661"
$ws.Range("D18").Value = "GLORIA ENIT DIAZ AVILA"
$ws.Range("E18").Value = "2111"
$ws.Range("F18").Value = 36341
$ws.Range("G18").Value = 781242

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "45761661"
$ws.Range("D19").Value = "GLORIA ENIT DIAZ AVILA"
$ws.Range("E19").Value = "2112"
$ws.Range("F19").Value = 36341
$ws.Range("G19").Value = 781242

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "45761661"
$ws.Range("D20").Value = "GLORIA ENIT DIAZ AVILA"
$ws.Range("E20").Value = "2204"
$ws.Range("F20").Value = 40000
$ws.Range("G20").Value = 781242

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "45761661"
$ws.Range("D21").Value = "GLORIA ENIT DIAZ AVILA"
$ws.Range("E21").Value = "2205"
$ws.Range("F21").Value = 30666
$ws.Range("G21").Value = 781242

# --- Column D is now narrower since the new name is shorter -----------------
$ws.Columns("D").ColumnWidth = 22.7265625
